$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove old hyperlinks (login creds replaced, no longer linked) ---
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A3").Hyperlinks.Delete()

# --- row 2: login sample data ---
$ws.Range("A2").Value = "arvknkh5m@mozmail.com"
$ws.Range("B2").Value = 123456
$ws.Range("C2").Font.Color = 2696481

# --- row 3: cleared out (old hyperlink row, keep its style only) ---
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""

# --- row 6 ---
$ws.Range("F6").Value = "test"

# --- row 8: "Search" banner ---
$ws.Range("A8").Value = "Search"
$ws.Range("A8:E8").HorizontalAlignment = -4108
$ws.Range("A8:E8").Merge()

# --- row 9-11: search fields ---
$ws.Range("A9").Value = "Search Key"
$ws.Range("A10").Value = "Cairo"
$ws.Range("A11").Value = ""

# --- row 13: "Add Listing" banner ---
$ws.Range("A13").Value = "Add Listing"
$ws.Range("A13:E13").HorizontalAlignment = -4108
$ws.Range("A13:E13").Merge()

# --- row 14: add-listing column headers ---
$ws.Range("A14").Value = "Email"
$ws.Range("B14").Value = "Password"
$ws.Range("C14").Value = "City"
$ws.Range("D14").Value = "District"
$ws.Range("E14").Value = "Street"
$ws.Range("F14").Value = "Furnishing status"
$ws.Range("G14").Value = "Property type"
$ws.Range("H14").Value = "Number of rooms"
$ws.Range("I14").Value = "Number of bathrooms"
$ws.Range("J14").Value = "Property name"
$ws.Range("K14").Value = "Rent / night"
$ws.Range("L14").Value = "Rent / month"
$ws.Range("M14").Value = "Description"

# --- row 15: add-listing sample data ---
$ws.Range("A15").Value = "autolandlord@1.com"
$ws.Range("B15").Value = 123456

# --- row 21: leftover test formula ---
$ws.Range("F21").Formula = "=0"

# --- column P: a vertical "." marker column running the whole used range (P10 stays empty) ---
$ws.Range("P1:P9").Value = "."
$ws.Range("P11:P27").Value = "."

# --- column widths / layout for the new sections ---
$ws.Columns("A").ColumnWidth = 17
$ws.Columns("C").ColumnWidth = 13.625
$ws.Columns("D").ColumnWidth = 12.875
$ws.Columns("F").ColumnWidth = 14.875
$ws.Columns("G").ColumnWidth = 11.625
$ws.Columns("H").ColumnWidth = 15
$ws.Columns("I").ColumnWidth = 18.75
$ws.Columns("J").ColumnWidth = 12.625

# --- selection to match new layout ---
$ws.Range("C15").Select()
